$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 7, 8, 9: the taxon/location records were re-ordered (a 3-way cyclic
# rotation of the per-observation fields: Id, Taxonsorteringsordning,
# TaxonId, Artnamn, Vetenskapligt namn, Ost, Nord, Publik kommentar).
#   new row7 <- old row8, new row8 <- old row9, new row9 <- old row7
# Only columns A,B,E,F,G,Q,R,AC actually change; the rest of each row
# (species author, dates, observer, etc.) stays put.
# ---------------------------------------------------------------------------

# --- Row 7 (becomes former row 8 data) ---
$ws.Range("A7").Value = 111815486
$ws.Range("B7").Value = 56398
$ws.Range("E7").Value = 100109
$ws.Range("F7").Value = "Tretåig hackspett"
$ws.Range("G7").Value = "Picoides tridactylus"
$ws.Range("Q7").Value = 457490.629824138
$ws.Range("R7").Value = 7057910.64054891
$ws.Range("AC7").Value = "ringhack"

# --- Row 8 (becomes former row 9 data) ---
$ws.Range("A8").Value = 111815489
$ws.Range("B8").Value = 56414
$ws.Range("E8").Value = 100049
$ws.Range("F8").Value = "Spillkråka"
$ws.Range("G8").Value = "Dryocopus martius"
$ws.Range("Q8").Value = 457851.1019836199
$ws.Range("R8").Value = 7058247.981310523
$ws.Range("AC8").Value = "hack"

# --- Row 9 (becomes former row 7 data) ---
$ws.Range("A9").Value = 111815490
$ws.Range("B9").Value = 56414
$ws.Range("Q9").Value = 457486.844484477
$ws.Range("R9").Value = 7058059.55768314

# ---------------------------------------------------------------------------
# Rows 22 and 23: the two observation records (Tretåig hackspett / bird, and
# Trådticka / fungus) were swapped entirely, including the extra
# Ålder-Stadium / Kön / Aktivitet / Metod (K:N) blank cells and the
# "Publik kommentar" (AC) comment, which only existed on the bird row.
# ---------------------------------------------------------------------------

# Move the K:N blank block and the AC comment from row 22 to row 23 first,
# so row 22 loses them and row 23 gains them (matching the swap).
$ws.Range("K22:N22").Cut($ws.Range("K23:N23"))
$ws.Range("AC22").Cut($ws.Range("AC23"))

# Now swap the remaining differing values between row 22 and row 23.
$ws.Range("A22").Value = 111815478
$ws.Range("B22").Value = 90087
$ws.Range("D22").Value = "LC"
$ws.Range("E22").Value = 3298
$ws.Range("F22").Value = "Trådticka"
$ws.Range("G22").Value = "Climacocystis borealis"
$ws.Range("H22").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("Q22").Value = 457490.7608241383
$ws.Range("R22").Value = 7057588.885967719

$ws.Range("A23").Value = 111815485
$ws.Range("B23").Value = 56398
$ws.Range("D23").Value = "NT"
$ws.Range("E23").Value = 100109
$ws.Range("F23").Value = "Tretåig hackspett"
$ws.Range("G23").Value = "Picoides tridactylus"
$ws.Range("H23").Value = "(Linnaeus, 1758)"
$ws.Range("Q23").Value = 457446.9368417656
$ws.Range("R23").Value = 7058136.079544679
